$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2073170731707317
$ws.Range("C2").Value = 0.5503048780487805
$ws.Range("J2").Value = 0.006097560975609756
$ws.Range("P2").Value = 0.1600609756097561
$ws.Range("S2").Value = 0.07621951219512195
$ws.Range("B3").Value = 0.002617801047120419
$ws.Range("C3").Value = 0.02879581151832461
$ws.Range("J3").Value = 0.02356020942408377
$ws.Range("P3").Value = 0.7643979057591623
$ws.Range("S3").Value = 0.1806282722513089
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.66
$ws.Range("S4").Value = 0.29
$ws.Range("J5").Value = 0.1428571428571428
$ws.Range("O5").Value = 0.1428571428571428
$ws.Range("P5").Value = 0.7142857142857143
$ws.Range("B6").Value = 0.06134969325153374
$ws.Range("D6").Value = 0.016359918200409
$ws.Range("E6").Value = 0.002044989775051125
$ws.Range("F6").Value = 0.06748466257668712
$ws.Range("J6").Value = 0.2310838445807771
$ws.Range("O6").Value = 0.01022494887525562
$ws.Range("Q6").Value = 0.1533742331288344
$ws.Range("R6").Value = 0.06748466257668712
$ws.Range("S6").Value = 0.3905930470347648
$ws.Range("B7").Value = 0.1113744075829384
$ws.Range("D7").Value = 0.01421800947867299
$ws.Range("F7").Value = 0.04976303317535545
$ws.Range("J7").Value = 0.1255924170616114
$ws.Range("O7").Value = 0.01658767772511848
$ws.Range("Q7").Value = 0.1753554502369668
$ws.Range("R7").Value = 0.08767772511848342
$ws.Range("S7").Value = 0.4194312796208531
$ws.Range("B8").Value = 0.1007268951194185
$ws.Range("D8").Value = 0.02388369678089304
$ws.Range("E8").Value = 0.003115264797507788
$ws.Range("F8").Value = 0.06022845275181724
$ws.Range("J8").Value = 0.122533748701973
$ws.Range("O8").Value = 0.01349948078920041
$ws.Range("Q8").Value = 0.1630321910695743
$ws.Range("R8").Value = 0.09138110072689512
$ws.Range("S8").Value = 0.4215991692627207
$ws.Range("B9").Value = 0.09828009828009827
$ws.Range("D9").Value = 0.01965601965601966
$ws.Range("F9").Value = 0.05651105651105651
$ws.Range("J9").Value = 0.1326781326781327
$ws.Range("O9").Value = 0.004914004914004914
$ws.Range("Q9").Value = 0.1597051597051597
$ws.Range("R9").Value = 0.09582309582309582
$ws.Range("S9").Value = 0.4324324324324325
$ws.Range("B10").Value = 0.1102592186929536
$ws.Range("D10").Value = 0.02154070828769624
$ws.Range("E10").Value = 0.001095290251916758
$ws.Range("F10").Value = 0.07228915662650602
$ws.Range("J10").Value = 0.1303395399780942
$ws.Range("O10").Value = 0.01387367652427893
$ws.Range("Q10").Value = 0.2066447608616283
$ws.Range("R10").Value = 0.0945600584154801
$ws.Range("S10").Value = 0.3493975903614458
$ws.Range("G11").Value = 0.1317957166392092
$ws.Range("J11").Value = 0.08896210873146623
$ws.Range("K11").Value = 0.1828665568369028
$ws.Range("L11").Value = 0.5815485996705108
$ws.Range("S11").Value = 0.01482701812191104
$ws.Range("F12").Value = 0.002710027100271003
$ws.Range("G12").Value = 0.7859078590785907
$ws.Range("J12").Value = 0.1490514905149052
$ws.Range("K12").Value = 0.008130081300813009
$ws.Range("L12").Value = 0.02710027100271003
$ws.Range("S12").Value = 0.02710027100271003
$ws.Range("F13").Value = 0.01111111111111111
$ws.Range("G13").Value = 0.7555555555555555
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.03333333333333333
$ws.Range("F15").Value = 0.02136752136752137
$ws.Range("H15").Value = 0.2158119658119658
$ws.Range("I15").Value = 0.06623931623931624
$ws.Range("J15").Value = 0.3632478632478632
$ws.Range("K15").Value = 0.06837606837606838
$ws.Range("M15").Value = 0.00641025641025641
$ws.Range("O15").Value = 0.07692307692307693
$ws.Range("S15").Value = 0.1816239316239316
$ws.Range("F16").Value = 0.02027027027027027
$ws.Range("H16").Value = 0.2027027027027027
$ws.Range("I16").Value = 0.06756756756756757
$ws.Range("J16").Value = 0.3986486486486486
$ws.Range("K16").Value = 0.1058558558558559
$ws.Range("M16").Value = 0.02252252252252252
$ws.Range("N16").Value = 0.004504504504504504
$ws.Range("O16").Value = 0.06531531531531531
$ws.Range("S16").Value = 0.1126126126126126
$ws.Range("F17").Value = 0.02152852529601722
$ws.Range("H17").Value = 0.1786867599569429
$ws.Range("I17").Value = 0.09364908503767493
$ws.Range("J17").Value = 0.4219590958019376
$ws.Range("K17").Value = 0.08826695371367062
$ws.Range("M17").Value = 0.02152852529601722
$ws.Range("N17").Value = 0.001076426264800861
$ws.Range("O17").Value = 0.077502691065662
$ws.Range("S17").Value = 0.09580193756727665
$ws.Range("F18").Value = 0.01535087719298246
$ws.Range("H18").Value = 0.1732456140350877
$ws.Range("I18").Value = 0.1030701754385965
$ws.Range("J18").Value = 0.3903508771929824
$ws.Range("K18").Value = 0.1030701754385965
$ws.Range("M18").Value = 0.01535087719298246
$ws.Range("O18").Value = 0.07894736842105263
$ws.Range("S18").Value = 0.1206140350877193
$ws.Range("F19").Value = 0.01722158438576349
$ws.Range("H19").Value = 0.2039800995024875
$ws.Range("I19").Value = 0.08151549942594719
$ws.Range("J19").Value = 0.3842326827401454
$ws.Range("K19").Value = 0.1056257175660161
$ws.Range("M19").Value = 0.02219670876387294
$ws.Range("N19").Value = 0.0003827018752391887
$ws.Range("O19").Value = 0.06926903941829315
$ws.Range("S19").Value = 0.115575966322235
